$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "SISTEMA" header in column F
$ws.Range("F1").Value = "SISTEMA"

# Copy formatting (style) from the previous header cell (E1) so the new
# header cell gets the same grey header style
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Set the width of the new column to match the template (stored width 17)
$ws.Columns.Item(6).ColumnWidth = 16.205607476635517

# Refresh the AutoFilter so it now spans A1:F1
$ws.AutoFilterMode = $false
$ws.Range("A1:F1").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# new AutoFilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$F`$1"
    }
}

# Select the new header cell, matching the resulting worksheet selection
$ws.Range("F1").Select()
